$d = $word.ActiveDocument

# The target is the last paragraph in the document body (an empty paragraph
# immediately preceding the sectPr), which currently has no runs. We need to
# append a run sequence to it, mixing Greek-language runs (carrying an
# explicit <w:rPr><w:lang w:val="el-GR"/></w:rPr>) with plain Latin-word runs
# (HTML / CSS / download / links) that carry no rPr at all.

$count = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($count)
$r = $p.Range

# Pull this paragraph's WordOpenXML (a full pkg:package wrapper) so we can
# splice in the new runs just before the paragraph's closing tag, then feed
# the modified package back in via InsertXML -- this is the only path that
# reliably round-trips explicit run-level <w:rPr><w:lang.../></w:rPr>
# formatting through this COM host (plain Range/Selection.LanguageID writes
# are not persisted for normal text ranges).
$openXml = $r.WordOpenXML

$newRuns =
  '<w:r><w:rPr><w:lang w:val="el-GR"/></w:rPr><w:t xml:space="preserve">Στο παραδοτέο πέρα του παρόντος αρχείου περιλαμβάνονται και τα 5 </w:t></w:r>' +
  '<w:r><w:t>HTML</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="el-GR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="el-GR"/></w:rPr><w:t xml:space="preserve">αρχεία, ένα </w:t></w:r>' +
  '<w:r><w:t>CSS</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="el-GR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="el-GR"/></w:rPr><w:t xml:space="preserve">αρχείο, μια εικόνα που υπάρχει στην αρχική σελίδα, καθώς και 4 ακόμα αρχεία τα οποία γίνονται </w:t></w:r>' +
  '<w:r><w:t>download</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="el-GR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="el-GR"/></w:rPr><w:t xml:space="preserve">από τα </w:t></w:r>' +
  '<w:r><w:t>links</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="el-GR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="el-GR"/></w:rPr><w:t>στις σελίδες ‘Έγγραφα Μαθήματος’ και ‘Εργασίες’</w:t></w:r>'

$marker = '</w:pPr></w:p>'
$idx = $openXml.IndexOf($marker)
if ($idx -lt 0) {
    throw "could not locate target paragraph end marker in WordOpenXML"
}
$insertPos = $idx + '</w:pPr>'.Length
$newOpenXml = $openXml.Substring(0, $insertPos) + $newRuns + $openXml.Substring($insertPos)

$null = $r.InsertXML($newOpenXml)
